# "bug with empty notes is resolved"
#
# 1. A footnote's text had picked up a stray trailing "a" (an artifact
#    left over from OCR/markup, e.g. "...པེ་ཅིན།a") that needs to be
#    dropped so the note reads "...པེ་ཅིན།".
# 2. The very last run of the body paragraph ("...བསྒྱུར་བའོ།།") needs
#    the page/folio marker "[༧༧ན]" appended to it.
# 3. The final footnote is an empty/orphaned note (just a stray "།"
#    punctuation mark, i.e. no real content) and both its reference in
#    the body and its entry in the footnotes part must be removed
#    entirely.
#
# NOTE: steps 1-2 are done *before* step 3 so the body insertion still
# lands on the run that is (at that point) the last one in the
# paragraph and merges into it, exactly mirroring what happens when a
# user edits text immediately before the footnote mark that follows it.

$d = $word.ActiveDocument

# --- 1. Fix the footnote whose text has a stray trailing "a" ---------
for ($i = 1; $i -le $d.Footnotes.Count; $i++) {
    $fn = $d.Footnotes.Item($i)
    $fnText = $fn.Range.Text
    if ($fnText.EndsWith("a")) {
        $fn.Range.Text = $fnText.Substring(0, $fnText.Length - 1)
    }
}

# --- 2. Append the "[༧༧ན]" folio marker to the colophon line ---------
$r = $d.Content
$found = $r.Find.Execute("ལོ་ཙཱ་བ་བློ་ལྡན་ཤེས་རབ་ཀྱིས་བསྒྱུར་བའོ།།", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $r.Collapse(0)
    $r.InsertAfter("[༧༧ན]")
}

# --- 3. Delete whichever footnote is now empty (no real content) -----
for ($i = $d.Footnotes.Count; $i -ge 1; $i--) {
    $fn = $d.Footnotes.Item($i)
    $bare = $fn.Range.Text -replace "[\s་།]", ""
    if ($bare.Length -eq 0) {
        $fn.Delete()
    }
}
